$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

function Set-PlainValue($cellRef, $value) {
    $ws.Range($cellRef).Value = $value
}

# Row 2
Set-PlainValue 'D2' '21.713.91'
Set-PlainValue 'E2' '  -1.47%  '

# Row 3
Set-PlainValue 'D3' '1.539.40'
Set-PlainValue 'E3' '  -1.00%  '

# Row 4
Set-PlainValue 'E4' '  +0.15%  '

# Row 5
Set-PlainValue 'E5' '  +0.03%  '

# Row 6
Set-TextValue 'D6' '289.83'
Set-PlainValue 'E6' '  +1.16%  '

# Row 7
Set-TextValue 'D7' '0.3895'
Set-PlainValue 'E7' '  +3.25%  '

# Row 8
Set-TextValue 'D8' '0.3168'
Set-PlainValue 'E8' '  -2.18%  '

# Row 9
Set-TextValue 'D9' '42.96'
Set-PlainValue 'E9' '  +4.31%  '

# Row 10
Set-TextValue 'D10' '0.07175'
Set-PlainValue 'E10' '  -1.86%  '

# Row 11
Set-TextValue 'D11' '1.054'
Set-PlainValue 'E11' '  -6.26%  '

# Row 12
Set-PlainValue 'E12' '  +0.09%  '

# Row 13
Set-TextValue 'D13' '5.613'
Set-PlainValue 'E13' '  -1.89%  '

# Row 14
Set-TextValue 'D14' '18.59'
Set-PlainValue 'E14' '  -4.00%  '

# Row 15
Set-TextValue 'D15' '6.610'
Set-PlainValue 'E15' '  -2.96%  '

# Row 16
Set-PlainValue 'D16' '1.544.99'
Set-PlainValue 'E16' '  -0.29%  '

# Row 17
Set-TextValue 'D17' '0.00001099'
Set-PlainValue 'E17' '  +1.69%  '

# Row 18
Set-TextValue 'D18' '0.06582'
Set-PlainValue 'E18' '  -1.07%  '

# Row 19
Set-TextValue 'D19' '83.10'
Set-PlainValue 'E19' '  -2.27%  '

# Row 20
Set-PlainValue 'E20' '  +0.02%  '

# Row 21
Set-TextValue 'D21' '6.143'
Set-PlainValue 'E21' '  -4.52%  '

# Row 22
Set-TextValue 'D22' '15.34'
Set-PlainValue 'E22' '  -4.07%  '

# Row 23
Set-TextValue 'D23' '10.87'
Set-PlainValue 'E23' '  -5.82%  '

# Row 24
Set-TextValue 'D24' '2.406'
Set-PlainValue 'E24' '  +7.30%  '

# Row 25
Set-PlainValue 'D25' '21.712.77'
Set-PlainValue 'E25' '  -1.49%  '

# Row 26
Set-TextValue 'D26' '2.357'
Set-PlainValue 'E26' '  -6.45%  '

# Row 27
Set-TextValue 'D27' '146.94'
Set-PlainValue 'E27' '  -2.08%  '

# Row 28
Set-TextValue 'D28' '18.35'
Set-PlainValue 'E28' '  -2.75%  '

# Row 29
Set-TextValue 'D29' '4.841'
Set-PlainValue 'E29' '  -0.14%  '

# Row 30
Set-PlainValue 'D30' '1.727.01'
Set-PlainValue 'E30' '  +0.01%  '

# Row 31
Set-TextValue 'D31' '117.37'
Set-PlainValue 'E31' '  -2.33%  '

# Row 32
Set-TextValue 'D32' '5.881'
Set-PlainValue 'E32' '  -0.57%  '

# Row 33
Set-TextValue 'D33' '0.9622'
Set-PlainValue 'E33' '  -14.75%  '

# Row 34
Set-TextValue 'D34' '0.08179'
Set-PlainValue 'E34' '  -0.06%  '

# Row 35
Set-TextValue 'D35' '8.802'
Set-PlainValue 'E35' '  -5.29%  '

# Row 36
Set-TextValue 'D36' '0.06071'
Set-PlainValue 'E36' '  -1.40%  '

# Row 37
Set-TextValue 'D37' '5.107'
Set-PlainValue 'E37' '  -2.36%  '

# Row 38
Set-TextValue 'D38' '0.02197'
Set-PlainValue 'E38' '  -3.88%  '

# Row 39
Set-TextValue 'D39' '0.2031'
Set-PlainValue 'E39' '  -3.88%  '

# Row 40
Set-PlainValue 'B40' 'WEMIXTOKEN'
Set-PlainValue 'C40' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D40' '1.432'
Set-PlainValue 'E40' '  -12.42%  '

# Row 41
Set-PlainValue 'B41' 'TrustWalletToken'
Set-PlainValue 'C41' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D41' '1.179'
Set-PlainValue 'E41' '  -2.97%  '

# Row 42
Set-PlainValue 'E42' '  -0.03%  '

# Row 43
Set-TextValue 'D43' '10.62'
Set-PlainValue 'E43' '  -2.57%  '

# Row 44
Set-TextValue 'D44' '0.5711'
Set-PlainValue 'E44' '  -3.94%  '

# Row 45
Set-TextValue 'D45' '3.736'
Set-PlainValue 'E45' '  +0.29%  '

# Row 46
Set-TextValue 'D46' '12.97'
Set-PlainValue 'E46' '  -4.29%  '

# Row 47
Set-TextValue 'D47' '0.5468'
Set-PlainValue 'E47' '  -4.65%  '

# Row 48
Set-PlainValue 'E48' '  +0.25%  '

# Row 49
Set-TextValue 'D49' '116.11'
Set-PlainValue 'E49' '  -3.14%  '

# Row 50
Set-TextValue 'D50' '1.865'
Set-PlainValue 'E50' '  -3.72%  '

# Row 51
Set-TextValue 'D51' '0.06701'
Set-PlainValue 'E51' '  -2.96%  '
